$d = $word.ActiveDocument

# 1. Paragraph 2 ("I am testing version control." + "  Ver. 1.1") already
#    reads as "I am testing version control.  Ver. 1.1" across two runs.
#    Use Find/Replace over the whole paragraph so the engine rewrites it
#    as a single run (matching the target OOXML) without touching the
#    neighbouring paragraphs.
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Find.Execute(
    "I am testing version control.  Ver. 1.1", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "I am testing version control.  Ver. 1.1", 2) | Out-Null

# 2. Append a brand-new paragraph at the end of the document: "Second
#    test of version control" -> Ver. 1.2, split across two runs the
#    same way the first "Ver. 1.1" paragraph was (base sentence run +
#    a trailing run holding just the incremented digit).
$endOfDoc = $d.Content
$endOfDoc.Collapse(0)
$newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>I am testing version control.  Ver. 1.</w:t></w:r><w:r><w:t>2</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endOfDoc.InsertXML($newParaXml) | Out-Null
